{"js": "const body = context.document.body;\nconst replacements = [\n  [\"18\u00d730=\", \"32\u00d764=\"],\n  [\"59\u00d774=\", \"58\u00d791=\"],\n  [\"22\u00d742=\", \"27\u00d757=\"],\n  [\"87\u00d730=\", \"47\u00d748=\"],\n  [\"83\u00d795=\", \"40\u00d717=\"],\n  [\"57\u00d789=\", \"37\u00d783=\"],\n  [\"97\u00d794=\", \"92\u00d780=\"],\n  [\"25\u00d751=\", \"50\u00d724=\"],\n  [\"79\u00d755=\", \"79\u00d766=\"],\n  [\"49\u00d777=\", \"51\u00d776=\"],\n  [\"66\u00d774=\", \"13\u00d730=\"],\n  [\"31\u00d769=\", \"30\u00d772=\"],\n  [\"30\u00d747=\", \"31\u00d773=\"],\n  [\"39\u00d795=\", \"96\u00d777=\"],\n  [\"99\u00d758=\", \"66\u00d755=\"],\n  [\"45\u00d759=\", \"75\u00d735=\"],\n  [\"81\u00d798=\", \"77\u00d720=\"],\n  [\"39\u00d736=\", \"62\u00d784=\"],\n  [\"65\u00d749=\", \"93\u00d766=\"],\n  [\"52\u00d750=\", \"14\u00d747=\"],\n  [\"67\u00d780=\", \"66\u00d751=\"],\n  [\"55\u00d784=\", \"74\u00d792=\"],\n  [\"11\u00d734=\", \"48\u00d778=\"],\n  [\"38\u00d799=\", \"91\u00d721=\"],\n  [\"65\u00d732=\", \"78\u00d797=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"18\u00d730=\"; New = \"32\u00d764=\" }\n  @{ Old = \"59\u00d774=\"; New = \"58\u00d791=\" }\n  @{ Old = \"22\u00d742=\"; New = \"27\u00d757=\" }\n  @{ Old = \"87\u00d730=\"; New = \"47\u00d748=\" }\n  @{ Old = \"83\u00d795=\"; New = \"40\u00d717=\" }\n  @{ Old = \"57\u00d789=\"; New = \"37\u00d783=\" }\n  @{ Old = \"97\u00d794=\"; New = \"92\u00d780=\" }\n  @{ Old = \"25\u00d751=\"; New = \"50\u00d724=\" }\n  @{ Old = \"79\u00d755=\"; New = \"79\u00d766=\" }\n  @{ Old = \"49\u00d777=\"; New = \"51\u00d776=\" }\n  @{ Old = \"66\u00d774=\"; New = \"13\u00d730=\" }\n  @{ Old = \"31\u00d769=\"; New = \"30\u00d772=\" }\n  @{ Old = \"30\u00d747=\"; New = \"31\u00d773=\" }\n  @{ Old = \"39\u00d795=\"; New = \"96\u00d777=\" }\n  @{ Old = \"99\u00d758=\"; New = \"66\u00d755=\" }\n  @{ Old = \"45\u00d759=\"; New = \"75\u00d735=\" }\n  @{ Old = \"81\u00d798=\"; New = \"77\u00d720=\" }\n  @{ Old = \"39\u00d736=\"; New = \"62\u00d784=\" }\n  @{ Old = \"65\u00d749=\"; New = \"93\u00d766=\" }\n  @{ Old = \"52\u00d750=\"; New = \"14\u00d747=\" }\n  @{ Old = \"67\u00d780=\"; New = \"66\u00d751=\" }\n  @{ Old = \"55\u00d784=\"; New = \"74\u00d792=\" }\n  @{ Old = \"11\u00d734=\"; New = \"48\u00d778=\" }\n  @{ Old = \"38\u00d799=\"; New = \"91\u00d721=\" }\n  @{ Old = \"65\u00d732=\"; New = \"78\u00d797=\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair.Old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair.New\n  $find.Execute(\n    [ref]$pair.Old,\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]$pair.New,\n    [ref]2\n  ) | Out-Null\n}\n"}
